$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A52").Value = 51
$ws.Range("B52").Value = 1
$ws.Range("C52").Value = "2024-06-15 23:13:02"
$ws.Range("D52").Value = 200
$ws.Range("E52").Value = 5

$ws.Range("A53").Value = 52
$ws.Range("B53").Value = 2
$ws.Range("C53").Value = "2024-06-15 23:13:02"
$ws.Range("D53").Value = 200
$ws.Range("E53").Value = 0
